$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.248.26"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.864.33"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "237.19"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "0.2867"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").Value = "0.06555"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "22.19"
$ws.Range("E10").Value = "  +10.09%  "
$ws.Range("D11").Value = "0.07892"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "97.97"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "1.869.42"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "5.197"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "0.6830"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "277.75"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "30.247.43"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "13.62"
$ws.Range("E18").Value = "  +7.92%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "0.000007351"
$ws.Range("D21").Value = "2.118.62"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "5.364"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "6.201"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Value = "168.19"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").Value = "9.253"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").Value = "19.08"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  +2.89%  "
$ws.Range("D29").Value = "1.383"
$ws.Range("E29").Value = "  +3.15%  "
$ws.Range("D30").Value = "0.09846"
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "1.482"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "4.076"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("D36").Value = "0.7054"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.01880"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "2.632"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "75.62"
$ws.Range("E40").Value = "  +4.08%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.289"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "1.958"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "0.8550"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "0.4182"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "103.58"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "7.222"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.290"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "949.46"
$ws.Range("E49").Value = "  -4.08%  "
$ws.Range("D50").Value = "34.27"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "0.05645"
$ws.Range("E51").Value = "  +0.15%  "
